# Generate Report for Handback
# The localization status workbook gets refreshed after a successful
# handback run: the overview "Status" column flips from "Ready for
# handoff" to "Handed back: in sync with en-US", and each language
# sheet (zh-cn / de-de) gets its "Latest Target File" / "Latest
# Handback File" columns populated (with hyperlinks) plus a fresh
# "Latest Handback DateTime" stamp.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1. Flip the "Ready for handoff" status to "Handed back: in sync
#    with en-US" everywhere it is used: the Overview rollup columns
#    (B/C) as well as the per-language "Status" column (C) on both
#    the zh-cn and de-de detail sheets.
# ---------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$newStatus = "Handed back: in sync with en-US"
$overview.Range("B2").Value = $newStatus
$overview.Range("C2").Value = $newStatus
$overview.Range("B3").Value = $newStatus
$overview.Range("C3").Value = $newStatus

# ---------------------------------------------------------------
# 2. zh-cn sheet: fill in Latest Target File (F) / Latest Handback
#    File (G) for both rows, and stamp the new handback datetime (H).
# ---------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("C3").Value = $newStatus

$zhcn.Hyperlinks.Add(
    $zhcn.Range("F2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/3e8012290e5153ac4fd467601a0878289068c5dd/e2e/ae0d94ec-2646-4100-8858-101a6a503f0e.md",
    "",
    "",
    "ae0d94ec-2646-4100-8858-101a6a503f0e.md"
)
$zhcn.Hyperlinks.Add(
    $zhcn.Range("G2"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8e3d2265248094417e6c3c769ab0f5df1a3b7be9/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/ae0d94ec-2646-4100-8858-101a6a503f0e.1e15a0e2b32aad7b1581d34b1d20e1277064190c.zh-cn.xlf",
    "",
    "",
    "ae0d94ec-2646-4100-8858-101a6a503f0e.1e15a0e2b32aad7b1581d34b1d20e1277064190c.zh-cn.xlf"
)
$zhcn.Hyperlinks.Add(
    $zhcn.Range("F3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/3e8012290e5153ac4fd467601a0878289068c5dd/e2e/f0067453-4126-4cb2-81ef-ade8d03f27c3.md",
    "",
    "",
    "f0067453-4126-4cb2-81ef-ade8d03f27c3.md"
)
$zhcn.Hyperlinks.Add(
    $zhcn.Range("G3"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8e3d2265248094417e6c3c769ab0f5df1a3b7be9/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/f0067453-4126-4cb2-81ef-ade8d03f27c3.c10901518470bf49267ab6ef4bddc88faba38a6f.zh-cn.xlf",
    "",
    "",
    "f0067453-4126-4cb2-81ef-ade8d03f27c3.c10901518470bf49267ab6ef4bddc88faba38a6f.zh-cn.xlf"
)

$zhcn.Range("H2").Value = "2016-03-12 08:33:24"
$zhcn.Range("H3").Value = "2016-03-12 08:33:24"

# ---------------------------------------------------------------
# 3. de-de sheet: same treatment, different timestamp/target files.
# ---------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

$dede.Hyperlinks.Add(
    $dede.Range("F2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/3e8012290e5153ac4fd467601a0878289068c5dd/e2e/ae0d94ec-2646-4100-8858-101a6a503f0e.md",
    "",
    "",
    "ae0d94ec-2646-4100-8858-101a6a503f0e.md"
)
$dede.Hyperlinks.Add(
    $dede.Range("G2"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d367b3246b093871854d2e3b3055c0f06c71c0b2/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/ae0d94ec-2646-4100-8858-101a6a503f0e.1e15a0e2b32aad7b1581d34b1d20e1277064190c.de-de.xlf",
    "",
    "",
    "ae0d94ec-2646-4100-8858-101a6a503f0e.1e15a0e2b32aad7b1581d34b1d20e1277064190c.de-de.xlf"
)
$dede.Hyperlinks.Add(
    $dede.Range("F3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/3e8012290e5153ac4fd467601a0878289068c5dd/e2e/f0067453-4126-4cb2-81ef-ade8d03f27c3.md",
    "",
    "",
    "f0067453-4126-4cb2-81ef-ade8d03f27c3.md"
)
$dede.Hyperlinks.Add(
    $dede.Range("G3"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d367b3246b093871854d2e3b3055c0f06c71c0b2/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/f0067453-4126-4cb2-81ef-ade8d03f27c3.c10901518470bf49267ab6ef4bddc88faba38a6f.de-de.xlf",
    "",
    "",
    "f0067453-4126-4cb2-81ef-ade8d03f27c3.c10901518470bf49267ab6ef4bddc88faba38a6f.de-de.xlf"
)

$dede.Range("H2").Value = "2016-03-12 08:33:30"
$dede.Range("H3").Value = "2016-03-12 08:33:30"
